{"js": "// Office.js (Word JavaScript API) edit script.\n// Goal (per commit message \"Added mention of creating T cell model\"):\n//   1. \"I also imported some animations\" -> \"I imported some animations\"\n//      in the paragraph about vertex arrays / models / animations.\n//   2. \"...for the player and UI elements.\" -> \"...for the player and UI\n//      elements as well as creating the model for T cells myself with\n//      blender. \" (note trailing space, matching the OOXML diff which adds\n//      a final run containing just a space).\n\nconst body = context.document.body;\n\n// --- Edit 1: \"I also imported some animations \" -> \"I imported some animations \" ---\nconst r1 = body.search(\"I also imported some animations\", { matchCase: true });\nr1.load(\"text\");\nawait context.sync();\n\nif (r1.items.length === 0) {\n  throw new Error(\"Could not find target text for edit 1 ('I also imported some animations').\");\n}\nr1.items[0].insertText(\"I imported some animations\", \"Replace\");\nawait context.sync();\n\n// --- Edit 2: append the T cell sentence after \"...UI elements.\" ---\nconst r2 = body.search(\n  \"for the player and UI elements.\",\n  { matchCase: true }\n);\nr2.load(\"text\");\nawait context.sync();\n\nif (r2.items.length === 0) {\n  throw new Error(\"Could not find target text for edit 2 ('for the player and UI elements.').\");\n}\nr2.items[0].insertText(\n  \"for the player and UI elements as well as creating the model for T cells myself with blender. \",\n  \"Replace\"\n);\nawait context.sync();\n", "ps1": "# Word COM interop script.\n# Goal (per commit message \"Added mention of creating T cell model\"):\n#   1. \"I also imported some animations\" -> \"I imported some animations\"\n#   2. \"...for the player and UI elements.\" -> \"...for the player and UI\n#      elements as well as creating the model for T cells myself with\n#      blender. \" (trailing space included, matching the OOXML diff).\n\n$d = $word.ActiveDocument\n\n# Word enumeration constants (not pre-seeded by this host, so spelled out):\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n# --- Edit 1: \"I also imported some animations\" -> \"I imported some animations\" ---\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$found1 = $find1.Execute(\n    \"I also imported some animations\",  # FindText\n    $true,                              # MatchCase\n    $false,                             # MatchWholeWord\n    $false,                             # MatchWildcards\n    $false,                             # MatchSoundsLike\n    $false,                             # MatchAllWordForms\n    $true,                              # Forward\n    $wdFindContinue,                    # Wrap\n    $false,                             # Format\n    \"I imported some animations\",       # ReplaceWith\n    $wdReplaceOne                       # Replace\n)\nif (-not $found1) {\n    throw \"Could not find target text for edit 1 ('I also imported some animations').\"\n}\n\n# --- Edit 2: append the T cell sentence after \"...UI elements.\" ---\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$found2 = $find2.Execute(\n    \"for the player and UI elements.\",  # FindText\n    $true,                              # MatchCase\n    $false,                             # MatchWholeWord\n    $false,                             # MatchWildcards\n    $false,                             # MatchSoundsLike\n    $false,                             # MatchAllWordForms\n    $true,                              # Forward\n    $wdFindContinue,                    # Wrap\n    $false,                             # Format\n    \"for the player and UI elements as well as creating the model for T cells myself with blender. \", # ReplaceWith\n    $wdReplaceOne                       # Replace\n)\nif (-not $found2) {\n    throw \"Could not find target text for edit 2 ('for the player and UI elements.').\"\n}\n"}
